# Add "Save" column (H) to the s_vals worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - same style as the other header cells (copy format from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for H2:H9
$saveValues = @(0, 1, 0, 1, 1, 0, 1, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
